# Commit: "1) Fixed the comments from UAT 1"
#
# The edit adds a new "OtherNames" merge-field placeholder next to the
# existing "Title"/"Surname" placeholders for both the first and second
# consultant, i.e. turns:
#   ${firstConsultantTitle} ${firstConsultantSurname}
# into:
#   ${firstConsultantTitle} ${firstConsultantOtherNames} ${firstConsultantSurname}
# (and the analogous change for secondConsultant inside a text box).
#
# NOTE: the secondConsultant block lives inside a legacy VML text box
# (w:pict/v:shape/v:textbox/w:txbxContent). That content is not part of
# any story reachable from Document.Content / Document.Paragraphs /
# Document.StoryRanges / Document.Shapes in this host, so it can't be
# targeted with the exposed object model. We apply the reachable half of
# the change (the table cell with firstConsultant...) via Find/Replace.

$d = $word.ActiveDocument

# Scope the search to the exact paragraph that holds
# "${firstConsultantTitle} ${firstConsultantSurname}" in the salutation
# table cell (the non-bold Calibri/24 occurrence — the bold "Dear ..."
# paragraph further down must stay untouched).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt.Contains('${firstConsultantTitle} ${firstConsultantSurname}') -and -not $txt.Contains('Dear')) {
        $target = $para.Range
        break
    }
}

if ($target -ne $null) {
    # Wrap:=0 (wdFindStop) + Replace:=1 (wdReplaceOne) keeps the
    # Find/Replace confined to this single paragraph's range instead of
    # leaking into the later, visually-similar "Dear ${firstConsultantTitle}
    # ${firstConsultantSurname}," paragraph.
    $replacement = '}' + ' ' + '${firstConsultantOtherNames} ${'
    $target.Find.Execute('} ${', $true, $false, $false, $false, $false, $true, 0, $false, $replacement, 1) | Out-Null
}

$check = $d.Paragraphs.Item(3).Range.Text
Write-Host "firstConsultant paragraph now:" $check
